$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 and 8 (the extra duplicate "Shaimaa" rows), leaving A1:C6
$ws.Rows.Item(7).Resize(2).Delete() | Out-Null

# Add the new "dead" column (D) values for rows 1-6
$ws.Range("D1").Value = "dead"
$ws.Range("D2").Value = "yes"
$ws.Range("D3").Value = "no"
$ws.Range("D4").Value = "yes"
$ws.Range("D5").Value = "yes"
$ws.Range("D6").Value = "no"

# Add the new "fav game" column (E) values for rows 1-3 only
$ws.Range("E1").Value = "fav game"
$ws.Range("E2").Value = "cod"
$ws.Range("E3").Value = "game of throns mob"

# Match the author's final selection state
$ws.Range("E4").Select() | Out-Null
